$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-12 (headers in row 1 stay the same: Area, Categoria, Proyecto, Lider, Video, Poster, Equipo)
$data = @(
    @("Nano",  "Prototipo",            "bgjonbipga´dnj´goa",                        "Marlon Yahir Martínez",  "FAS",          "FAS",         "bgjonbipga´dnj´goa"),
    @("Nano",  "Prototipo finalizado", "ok",                                         "Marlon Yahir Martínez",  "fsa",          "fas",         "ok"),
    @("Nano",  "Prototipo finalizado", "fsa",                                        "Marlon Yahir Martínez",  "fsaf",         "fsa",         "fsa"),
    @("Nano",  "Prototipo finalizado", "Proyecto para pbrar el commit de gerry",     "Marlon Yahir Martínez",  "fsfa",         "fasfa",       "Proyecto para pbrar el commit de gerry"),
    @("Bio",   "Prototipo",            "Otra nueva actualización",                   "Marlon Yahir Martínez",  "video123",     "poster123",   "Proyecto de prueba para autorizar"),
    @("Bio",   "Prototipo",            "Prueba de que funciona",                     "Marlon Yahir Martínez",  "fasfsa",       "fsafas",      "Prueba de que funciona"),
    @("Nano",  "Prototipo",            "Proyecto de prueba",                         "Marlon Yahir Martínez",  "video.pdf",    "poster.pdf",  "Proyecto de prueba"),
    @("Nexus", "Prototipo",            "Esto funcionara ",                           "Marlon Yahir Martínez",  "fasfa",        "fsfafsa",     "Esto funcionara "),
    @("Nano",  "Prototipo finalizado", "fad",                                        "Marlon Yahir Martínez",  "fasfaf",       "fsafa",       "fad"),
    @("Nexus", "Concepto",             "Prueba para certificados no lo borren",      "Sarai Santiago Lozano",  "holaaaaaaaa2", "holaaaaa",    "Prueba para certificados no lo borren"),
    @("Nexus", "Prototipo",            "Robot automata para automatizar automatas",  "Gerry Deustúa Hernández","DBDFBDFB",     "HDD",         "BOLDBGOSDBGOSANGSIGNSOPGSDG")
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}

# Copy the formatting of the existing data row (row 3) down across the newly added rows (4-12)
$ws.Range("A3:G3").Copy()
$ws.Range("A4:G12").PasteSpecial(-4122)

